$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.327.93"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "2.178.77"
$ws.Range("E3").Value = "  -1.65%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'252.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.60%  "
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("D7").Value = "'72.89"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.48%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'0.580"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.14%  "
$ws.Range("D10").Value = "'40.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.33%  "
$ws.Range("D11").Value = "'0.0907"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.18%  "
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("D13").Value = "'6.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.55%  "
$ws.Range("D14").Value = "2.504.87"
$ws.Range("E14").Value = "  -1.57%  "
$ws.Range("D15").Value = "'14.15"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.90%  "
$ws.Range("D16").Value = "2.179.60"
$ws.Range("E16").Value = "  -1.51%  "
$ws.Range("D17").Value = "'0.767"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.52%  "
$ws.Range("D18").Value = "42.264.17"
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("E19").Value = "  -2.95%  "
$ws.Range("D20").Value = "'70.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("E21").Value = "  -1.78%  "
$ws.Range("D22").Value = "'226.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.61%  "
$ws.Range("D23").Value = "'9.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.95%  "
$ws.Range("D24").Value = "'2.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.42%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("E26").Value = "  -4.49%  "
$ws.Range("D27").Value = "'3.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.26%  "
$ws.Range("E28").Value = "  -2.26%  "
$ws.Range("D29").Value = "'2.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("D30").Value = "'36.65"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.70%  "
$ws.Range("D31").Value = "'169.99"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.68%  "
$ws.Range("D32").Value = "'19.96"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.68%  "
$ws.Range("D33").Value = "'0.0816"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.30%  "
$ws.Range("D34").Value = "'5.09"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.62%  "
$ws.Range("D36").Value = "'0.106"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.99%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "'4.18"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.52%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.0333"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.96%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").Value = "'2.04"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.11%  "
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").Value = "'11.70"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.07%  "
$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D41").Value = "'59.13"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.36%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.194"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.90%  "
$ws.Range("D43").Value = "'5.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.83%  "
$ws.Range("D44").Value = "'101.53"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.26%  "
$ws.Range("E45").Value = "  +6.91%  "
$ws.Range("D46").Value = "'0.0968"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.21%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'8.14"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.67%  "
$ws.Range("B48").Value = "WOONetwork"
$ws.Range("C48").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D48").Value = "'0.457"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.44%  "
$ws.Range("E49").Value = "  -2.57%  "
$ws.Range("E50").Value = "  -1.46%  "
$ws.Range("E51").Value = "  +0.27%  "
